# Updated files for latest order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("N1").Value = "Order QTY"
$ws.Range("O1").Value = "Price"

# Row 2 (own formula, not part of the shared group below)
$ws.Range("N2").Formula = "=J2/10"
$ws.Range("O2").Value = 0.127

# Rows 3-12 share one formula (Excel fills it as a shared formula group)
$ws.Range("N3:N12").Formula = "=J3/10"
# Row 10 has no data in this table, so drop the filled cell there
$ws.Range("N10").ClearContents()

$ws.Range("O3").Value = 0.094
$ws.Range("O4").Value = 0.01
$ws.Range("O5").Value = 0.042
$ws.Range("O6").Value = 0.08
$ws.Range("O7").Value = 0.078
$ws.Range("O8").Value = 0.669
$ws.Range("O9").Value = 0.327
$ws.Range("O11").Value = 0.244
$ws.Range("O12").Value = 0.96

# Totals row
$ws.Range("K14").Formula = "=SUM(K2:K12)"
$ws.Range("K14").Style = "Normal"
$ws.Range("O14").Formula = "=SUM(O2:O12)"

# Formatting-only cell (currency style, no value)
$ws.Range("K16").NumberFormat = $ws.Range("K2").NumberFormat

# Update view/selection
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J14").Select()
